# The post "「20%学生割引」" (row 131) was removed from the posts list.
# Deleting the entire row shifts every row below it up by one, which is
# exactly what the target diff shows (row 132 -> 131, ..., row 328 -> 327)
# and also updates the sheet's dimension from A1:C328 to A1:C327.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(131).Delete()
